$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 79; existing rows 79:85 shift down to 81:87.
$ws.Rows.Item(79).Resize(2).Insert()

# New row 79 values
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(79, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(79, 4).Value = 45212
$ws.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = 100112027
$ws.Cells.Item(79, 7).Value = "Melón"
$ws.Cells.Item(79, 8).Value = "Tuna"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 50
$ws.Cells.Item(79, 11).Value = 24000
$ws.Cells.Item(79, 12).Value = 25000
$ws.Cells.Item(79, 13).Value = 24500
$ws.Cells.Item(79, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 1361
$ws.Cells.Item(79, 17).Value = 18
$ws.Cells.Item(79, 18).Value = "Hortaliza"

# New row 80 values
$ws.Cells.Item(80, 1).Value = 1
$ws.Cells.Item(80, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(80, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(80, 4).Value = 45212
$ws.Cells.Item(80, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(80, 5).Value = 15
$ws.Cells.Item(80, 6).Value = 100112027
$ws.Cells.Item(80, 7).Value = "Melón"
$ws.Cells.Item(80, 8).Value = "Tuna"
$ws.Cells.Item(80, 9).Value = "Segunda"
$ws.Cells.Item(80, 10).Value = 50
$ws.Cells.Item(80, 11).Value = 22000
$ws.Cells.Item(80, 12).Value = 23000
$ws.Cells.Item(80, 13).Value = 22540
$ws.Cells.Item(80, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(80, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(80, 16).Value = 939
$ws.Cells.Item(80, 17).Value = 24
$ws.Cells.Item(80, 18).Value = "Hortaliza"
